# Add initial conditions for bioreactors to the "asm1" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("asm1")
$wwSheet = $wb.Worksheets.Item("wastewater")

# Header row (row 1) - matches the headers already used on the "wastewater" sheet
$ws.Range("B1").Value = "S_I"
$ws.Range("C1").Value = "S_S"
$ws.Range("D1").Value = "X_I"
$ws.Range("E1").Value = "X_S"
$ws.Range("F1").Value = "X_BH"
$ws.Range("G1").Value = "X_BA"
$ws.Range("H1").Value = "X_P"
$ws.Range("I1").Value = "S_O"
$ws.Range("J1").Value = "S_NO"
$ws.Range("K1").Value = "S_NH"
$ws.Range("L1").Value = "S_ND"
$ws.Range("M1").Value = "X_ND"
$ws.Range("N1").Value = "S_ALK"

# Row labels
$ws.Range("A2").Value = "A1"
$ws.Range("A3").Value = "A2"
$ws.Range("A4").Value = "O1"
$ws.Range("A5").Value = "O2"
$ws.Range("A6").Value = "O3"
$ws.Range("A7").Value = "C1_s"
$ws.Range("A8").Value = "C1_x"
$ws.Range("A9").Value = "C1_tss"

# Row 2 data (A1)
$ws.Range("B2").Value = 28.0643
$ws.Range("C2").Value = 3.0503
$ws.Range("D2").Value = 1532.3
$ws.Range("E2").Value = 63.0433
$ws.Range("F2").Value = 2245.1
$ws.Range("G2").Value = 166.6699
$ws.Range("H2").Value = 964.8992
$ws.Range("I2").Value = 0.0093
$ws.Range("J2").Value = 3.935
$ws.Range("K2").Value = 6.8924
$ws.Range("L2").Value = 0.958
$ws.Range("M2").Value = 3.8453
$ws.Range("N2").Formula = "=5.4213*12"

# Row 3 data (A2)
$ws.Range("B3").Value = 28.0643
$ws.Range("C3").Value = 1.3412
$ws.Range("D3").Value = 1532.3
$ws.Range("E3").Value = 58.8579
$ws.Range("F3").Value = 2245.4
$ws.Range("G3").Value = 166.5512
$ws.Range("H3").Value = 965.6805
$ws.Range("I3").Value = 0.00010907
$ws.Range("J3").Value = 2.2207
$ws.Range("K3").Value = 7.2028
$ws.Range("L3").Value = 0.6862
$ws.Range("M3").Value = 3.7424
$ws.Range("N3").Formula = "=5.5659*12"

# Row 4 data (O1)
$ws.Range("B4").Value = 28.0643
$ws.Range("C4").Value = 0.9553
$ws.Range("D4").Value = 1532.3
$ws.Range("E4").Value = 46.2983
$ws.Range("F4").Value = 2246.8
$ws.Range("G4").Value = 167.3077
$ws.Range("H4").Value = 967.2442
$ws.Range("I4").Value = 0.4663
$ws.Range("J4").Value = 5.5141
$ws.Range("K4").Value = 3.4247
$ws.Range("L4").Value = 0.6513
$ws.Range("M4").Value = 3.1405
$ws.Range("N4").Formula = "=5.0608*12"

# Row 5 data (O2)
$ws.Range("B5").Value = 28.0643
$ws.Range("C5").Value = 0.7806
$ws.Range("D5").Value = 1532.3
$ws.Range("E5").Value = 37.3881
$ws.Range("F5").Value = 2245.6
$ws.Range("G5").Value = 167.8339
$ws.Range("H5").Value = 968.8072
$ws.Range("I5").Value = 1.4284
$ws.Range("J5").Value = 8.4066
$ws.Range("K5").Value = 0.6922
$ws.Range("L5").Value = 0.6094
$ws.Range("M5").Value = 2.6815
$ws.Range("N5").Formula = "=4.659*12"

# Row 6 data (O3)
$ws.Range("B6").Value = 28.0643
$ws.Range("C6").Value = 0.6734
$ws.Range("D6").Value = 1532.3
$ws.Range("E6").Value = 31.9144
$ws.Range("F6").Value = 2242.1
$ws.Range("G6").Value = 167.8482
$ws.Range("H6").Value = 970.3678
$ws.Range("I6").Value = 1.3748
$ws.Range("J6").Value = 9.1948
$ws.Range("K6").Value = 0.1585
$ws.Range("L6").Value = 0.5594
$ws.Range("M6").Value = 2.3926
$ws.Range("N6").Formula = "=4.5646*12"

$ws.Range("A9").Select()

# Residual conditional-format swatches (Light Red Fill with Dark Red Text)
# left behind in the style pool, matching the workbook's saved state.
$cfRange = $wwSheet.Range("N2:N3")

$fc1 = $cfRange.FormatConditions.AddUniqueValues()
$fc1.DupeUnique = 1
$fc1.Font.Color = 393372
$fc1.Interior.Color = 13551615
$fc1.Delete()

$fc2 = $cfRange.FormatConditions.AddUniqueValues()
$fc2.DupeUnique = 1
$fc2.Font.Color = 393372
$fc2.Interior.Color = 13551615
$fc2.Delete()
